$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# repull data, push all data, mean calculation
# Update dSF (column F) values for several rows to reflect repulled data
$ws.Range("F16").Value = 1
$ws.Range("F25").Value = 2
$ws.Range("F26").Value = 2
$ws.Range("F28").Value = -3
